$wb = $excel.ActiveWorkbook

# Both the "展览" sheet and the "全部类型" sheet contain the same table and
# need their F2:F5 ("想去人数") values updated identically.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 6892
    $ws.Range("F3").Value = 52
    $ws.Range("F4").Value = 201
    $ws.Range("F5").Value = 37
}
